# Rename the 14 "(1SET)" / "(1PC)" suffixed item names in A45:A58 by
# stripping the trailing qualifier, matching the author's upload edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write order matches the order new shared-string entries were appended in
# the authored workbook (not simple row order).
$ws.Range("A46").Value = "W3 #6 機 改造 驅動輥(1SET=4PCS) 舊品拆裝"
$ws.Range("A47").Value = "W3X 串列式驅動輥軸心校正、研磨"
$ws.Range("A49").Value = "W3 #6 機 改造 驅動輥內孔車修(1SET=4PCS) 內孔車修"
$ws.Range("A48").Value = "W3 #6 機 改造 驅動輥內孔磨修(1SET=4PCS) 內孔手工研磨"
$ws.Range("A50").Value = "W3 #6 機 改造 240 驅動輥 (1SET=4PCS) 本體銲補"
$ws.Range("A51").Value = "W3 #6 機 改造 265 驅動輥 (1SET=4PCS) 本體銲補"
$ws.Range("A52").Value = "W3 #6 機 改造 295 驅動輥 (1SET=4PCS) 本體銲補"
$ws.Range("A53").Value = "W3 #6 機 改造 240 驅動輥 (1SET=4PCS) 本體未再生車修"
$ws.Range("A54").Value = "W3 #6 機 改造 265 驅動輥 (1SET=4PCS) 本體未再生車修"
$ws.Range("A56").Value = "W3 #6 機 改造 240 驅動輥 (1SET=4PCS) 本體再生車修"
$ws.Range("A55").Value = "W3 #6 機 改造 295 驅動輥 (1SET=4PCS) 本體未再生車修"
$ws.Range("A58").Value = "W3 #6 機 改造 295 驅動輥 (1SET=4PCS) 本體再生車修"
$ws.Range("A57").Value = "W3 #6 機 改造 265 驅動輥 (1SET=4PCS) 本體再生車修"
$ws.Range("A45").Value = "W3 #6 機 改造 驅動輥(1SET=4PCS) 新品組裝"

# Match the saved cursor/selection position from the authored workbook.
$ws.Range("A49").Select()
